# Nick Tackes CV — apply commit "publications and diss fellowship"
#
# 1. Update the CV date 2021-06-09 -> 2021-06-29 (title-page date paragraph
#    and the cached DATE field result in the running header).
# 2. Add a new "2021 | Columbia University IRCPL Dissertation Fellowship"
#    row to the "Competitive Scholarships and Honors" table, right before
#    the existing 2018 row.
# 3. Light proofreading touch-ups in the Publications / Conferences /
#    Teaching sections: re-save the (unchanged) paragraph text so that
#    Word re-normalizes run boundaries, matching the author's pass over
#    those paragraphs.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1a. Title-page date (plain run of text in the body).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("2021-06-09", $false, $false, $false, $false, `
    $false, $true, 1, $false, "2021-06-29", 2) | Out-Null

# ---------------------------------------------------------------------
# 1b. Header date field (a live DATE field with a cached result). The
#     cached text is only reliably reachable character-by-character (the
#     field's live Result.Text recomputes "today", but Result.Characters
#     still exposes the saved cache), so the ten characters of
#     "2021-06-09" are checked one at a time and only the single
#     character that actually differs ('0' -> '2') is written back. That
#     keeps the begin/instrText/separate/end run structure intact instead
#     of collapsing/reordering it the way a whole-range Find/Replace does.
# ---------------------------------------------------------------------
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)
    for ($h = 1; $h -le $section.Headers.Count; $h++) {
        $header = $section.Headers.Item($h)
        for ($fi = 1; $fi -le $header.Range.Fields.Count; $fi++) {
            $fld = $header.Range.Fields.Item($fi)
            $chars = $fld.Result.Characters
            if ($chars.Count -ge 10) {
                $c1 = $chars.Item(1).Text
                $c2 = $chars.Item(2).Text
                $c3 = $chars.Item(3).Text
                $c4 = $chars.Item(4).Text
                $c5 = $chars.Item(5).Text
                $c6 = $chars.Item(6).Text
                $c7 = $chars.Item(7).Text
                $c8 = $chars.Item(8).Text
                $c9 = $chars.Item(9).Text
                $c10 = $chars.Item(10).Text
                $isDate = ($c1 -eq "2") -and ($c2 -eq "0") -and ($c3 -eq "2") `
                    -and ($c4 -eq "1") -and ($c5 -eq "-") -and ($c6 -eq "0") `
                    -and ($c7 -eq "6") -and ($c8 -eq "-") -and ($c9 -eq "0") `
                    -and ($c10 -eq "9")
                if ($isDate) {
                    $chars.Item(9).Text = "2"
                }
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. New fellowship row in the "Competitive Scholarships and Honors"
#    table — inserted right before the row that currently starts with
#    "2018".
# ---------------------------------------------------------------------
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables.Item($ti)
    if ($tbl.Cell(1, 1).Range.Text.TrimEnd([char]13, [char]7) -eq "2018") {
        $newRow = $tbl.Rows.Add($tbl.Rows.Item(1))
        $newRow.Cells.Item(1).Range.Text = "2021"
        $newRow.Cells.Item(2).Range.Text = "Columbia University IRCPL Dissertation Fellowship"
        break
    }
}

# ---------------------------------------------------------------------
# 3. Re-apply the (unchanged) text of a handful of paragraphs that the
#    author revisited while proofreading the Publications / Conferences /
#    Teaching sections. The content is identical; this only normalizes
#    run boundaries the way Word does on a real edit pass.
# ---------------------------------------------------------------------
function Touch-Text([string]$text) {
    $rng = $word.ActiveDocument.Content
    $rng.Find.Execute($text, $false, $false, $false, $false, $false, `
        $true, 1, $false, $text, 2) | Out-Null
}

Touch-Text "“Om Shanti Emojis: Three Facets of Digital Hinduism,” Anthropology of Religion Unit and Religion, Media, and Culture Unit, American Academy of Religion, Online, December 5, 2020."
Touch-Text "“Energy and Vibrations: The Logic of Transformation in the Gayatri Pariwar and the Brahma Kumaris,” Public Health Workshop, Jalaharwal Nehru University, New Delhi, March 12, 2020."
Touch-Text "“Zooming in on Mozoomdar: A Microhistory of Brahmo Belief,” Religion in South Asia Section, American Academy of Religion, Denver, November 18, 2018."
Touch-Text "Chair/Discussant, “Yoga and Politics: South Asia and Beyond,” Madison South Asia Conference, October 12, 2018."
Touch-Text "“COVID-19 First Responders: The Gayatri Pariwar and the Immune Ritual Body.” "
Touch-Text "“Metabolic Living: Food, Fat, and the Absorption of Illness in India by Harris Solomon (Review).” "
Touch-Text "Lead Teaching Fellow (Center for Teaching and Learning, Columbia University, 2020-2021)"
Touch-Text "“East Asian Buddhism” (Michael Como, Religion)"
